$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.824.61"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.908.64"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.41"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4990"
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3786"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07272"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.19"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9031"
$ws.Range("E11").Value = "  -3.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07638"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "1.895.07"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.467"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.10"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008719"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "27.845.72"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.58"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.168"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "2.127.54"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.84"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.601"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.15"
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.846"
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.216"
$ws.Range("E27").Value = "  +4.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.37"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.15"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.872"
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08973"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.194"
$ws.Range("E32").Value = "  -1.97%  "

# Rows 33 and 34 swap places: ImmutableX <-> ARBITRUM
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.238"
$ws.Range("E33").Value = "  -1.28%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7888"
$ws.Range("E34").Value = "  +2.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.803"
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.652"
$ws.Range("E36").Value = "  +2.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02078"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.057"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5516"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05294"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.770"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.05"
$ws.Range("E43").Value = "  +3.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.468"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.54"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4793"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.634"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.32"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06029"
$ws.Range("E51").Value = "  -0.76%  "
